$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 with Jane's data
$ws.Range("A3").Value = "Jane"

# Update row 2: Last Name and Email change
$ws.Range("B2").Value = "Smith"
$ws.Range("C2").Value = "john.smith@fake.com"

# Continue filling in row 3
$ws.Range("B3").Value = "Doe"
$ws.Range("C3").Value = "jane.doe@fake.com"
$ws.Range("D3").Value = "123-123-1234"
$ws.Range("E3").Value = "ACC124"
$ws.Range("F3").Value = "Adult"
$ws.Range("G3").Value = "No"
$ws.Range("H3").Value = "Yes"
$ws.Range("I3").Value = 34.99
$ws.Range("I3").NumberFormat = $ws.Range("I2").NumberFormat
$ws.Range("J3").Value = "No"
$ws.Range("L3").Value = "Music is too loud."

# Column widths adjust (best fit) to accommodate the new/longer values.
$ws.Columns.Item(3).ColumnWidth = 18.5
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666

$ws.Range("F14").Select()
